# Recomputed TPM-based NATMI metrics for Bmp6 -> Bmpr1a (ligand/receptor detection,
# expression, specificity and edge-weight columns E:T) for all 9 cluster pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 57.14035266666667
$ws.Range("H2").Value = 171.421058
$ws.Range("I2").Value = 0.7274038390747541
$ws.Range("J2").Value = 0.7274038390747541
$ws.Range("M2").Value = 2.341355666666667
$ws.Range("N2").Value = 7.024067000000001
$ws.Range("O2").Value = 0.03973512964576821
$ws.Range("P2").Value = 0.0397351296457682
$ws.Range("Q2").Value = 133.7858885114318
$ws.Range("R2").Value = 1204.072996602886
$ws.Range("S2").Value = 0.02890348585046487
$ws.Range("T2").Value = 0.02890348585046487

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 57.14035266666667
$ws.Range("H3").Value = 171.421058
$ws.Range("I3").Value = 0.7274038390747541
$ws.Range("J3").Value = 0.7274038390747541
$ws.Range("O3").Value = 0.5779093692199981
$ws.Range("P3").Value = 0.5779093692199981
$ws.Range("Q3").Value = 1945.787496591513
$ws.Range("R3").Value = 17512.08746932362
$ws.Range("S3").Value = 0.4203734938078962
$ws.Range("T3").Value = 0.4203734938078962

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 57.14035266666667
$ws.Range("H4").Value = 171.421058
$ws.Range("I4").Value = 0.7274038390747541
$ws.Range("J4").Value = 0.7274038390747541
$ws.Range("O4").Value = 0.3823555011342337
$ws.Range("P4").Value = 0.3823555011342337
$ws.Range("Q4").Value = 1287.368907626682
$ws.Range("R4").Value = 11586.32016864014
$ws.Range("S4").Value = 0.2781268594163931
$ws.Range("T4").Value = 0.2781268594163931

# Row 5
$ws.Range("I5").Value = 0.08622113322131104
$ws.Range("J5").Value = 0.08622113322131104
$ws.Range("M5").Value = 2.341355666666667
$ws.Range("N5").Value = 7.024067000000001
$ws.Range("O5").Value = 0.03973512964576821
$ws.Range("P5").Value = 0.0397351296457682
$ws.Range("Q5").Value = 15.85800114988144
$ws.Range("R5").Value = 142.722010348933
$ws.Range("S5").Value = 0.003426007906753847
$ws.Range("T5").Value = 0.003426007906753846

# Row 6
$ws.Range("I6").Value = 0.08622113322131104
$ws.Range("J6").Value = 0.08622113322131104
$ws.Range("O6").Value = 0.5779093692199981
$ws.Range("P6").Value = 0.5779093692199981
$ws.Range("S6").Value = 0.04982800071336129
$ws.Range("T6").Value = 0.04982800071336129

# Row 7
$ws.Range("I7").Value = 0.08622113322131104
$ws.Range("J7").Value = 0.08622113322131104
$ws.Range("O7").Value = 0.3823555011342337
$ws.Range("P7").Value = 0.3823555011342337
$ws.Range("S7").Value = 0.03296712460119591
$ws.Range("T7").Value = 0.03296712460119591

# Row 8
$ws.Range("I8").Value = 0.1863750277039348
$ws.Range("J8").Value = 0.1863750277039348
$ws.Range("M8").Value = 2.341355666666667
$ws.Range("N8").Value = 7.024067000000001
$ws.Range("O8").Value = 0.03973512964576821
$ws.Range("P8").Value = 0.0397351296457682
$ws.Range("Q8").Value = 34.278549738519
$ws.Range("R8").Value = 308.5069476466711
$ws.Range("S8").Value = 0.00740563588854949
$ws.Range("T8").Value = 0.00740563588854949

# Row 9
$ws.Range("I9").Value = 0.1863750277039348
$ws.Range("J9").Value = 0.1863750277039348
$ws.Range("O9").Value = 0.5779093692199981
$ws.Range("P9").Value = 0.5779093692199981
$ws.Range("Q9").Value = 498.548645336397
$ws.Range("R9").Value = 4486.937808027573
$ws.Range("S9").Value = 0.1077078746987406
$ws.Range("T9").Value = 0.1077078746987406

# Row 10
$ws.Range("I10").Value = 0.1863750277039348
$ws.Range("J10").Value = 0.1863750277039348
$ws.Range("O10").Value = 0.3823555011342337
$ws.Range("P10").Value = 0.3823555011342337
$ws.Range("S10").Value = 0.07126151711664466
$ws.Range("T10").Value = 0.07126151711664468
